# Updates the crypto price/volume table (and the Stellar/Toncoin row swap)
# per the "Updated cryptos list" GitHub Actions commit.
#
# D-column cells whose new text parses as a plain number ("246.33",
# "0.9996", ...) are written with NumberFormat "@" (Text) first so Excel
# stores the literal string instead of silently coercing it to a number;
# the format is then reset to "Normal" so no stray style survives on the
# cell (matching the original file, which has no per-cell style there).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.421.16"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.869.61"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4741"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2919"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +6.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7388"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.19%  "
$ws.Range("D14").Value = "1.871.62"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.126"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "30.408.89"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007540"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.118.91"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.228"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.169"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.929"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.370"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09988"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.503"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.304"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.140"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04838"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01860"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.741"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.966"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4190"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8337"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.223"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "931.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05641"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
